$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test results / status / dates (simulating a fresh test run writing
# its results - e.g. a "Create Header in doc and Set value in header" run).
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "Failed"
$ws.Range("H2").Value = "30/05/2020"

$ws.Range("B3").Value = "No"
$ws.Range("H3").Value = "26/05/2020"

$ws.Range("B5").Value = "Yes"
$ws.Range("H5").Value = "30/05/2020"

# Move the active selection to B6, matching the sheet view after the run.
$ws.Range("B6").Select()
